$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# Typography sheet: "Default" typography's Wildcard Characters changes from "0-9" to "."
$wsTypo.Range("G4").Value = "."

# Translation sheet: ADC value text now includes a trailing "v" (volts)
$wsTrans.Range("F4").Value = "ADC value = <value>v"

# Translation sheet: example text changes from "1000" to "3.301" (an ADC voltage reading).
# Prefix with an apostrophe so the numeric-looking text stays plain text, not a number.
$wsTrans.Range("F5").Value = "'3.301"
